$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.691.39'
$ws.Range('E2').Value = '  +0.85%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.302.97'
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.65'
$ws.Range('E5').Value = '  -0.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.83'
$ws.Range('E6').Value = '  -0.44%  '
$ws.Range('E7').Value = '  -0.38%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.605'
$ws.Range('E9').Value = '  -0.69%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.84'
$ws.Range('E10').Value = '  +0.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0901'
$ws.Range('E11').Value = '  -0.84%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.51'
$ws.Range('E12').Value = '  +2.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.997'
$ws.Range('E14').Value = '  +3.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.40'
$ws.Range('E15').Value = '  +0.63%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.651.99'
$ws.Range('E16').Value = '  +0.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.300.70'
$ws.Range('E17').Value = '  +0.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.627.62'
$ws.Range('E18').Value = '  +1.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.61'
$ws.Range('E19').Value = '  +3.43%  '
$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.91'
$ws.Range('E20').Value = '  +30.86%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000106'
$ws.Range('E21').Value = '  +0.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '74.06'
$ws.Range('E22').Value = '  +0.98%  '
$ws.Range('E23').Value = '  -1.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '266.66'
$ws.Range('E24').Value = '  -4.03%  '
$ws.Range('E25').Value = '  -1.31%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('E27').Value = '  +0.66%  '
$ws.Range('E28').Value = '  -4.30%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '22.64'
$ws.Range('E29').Value = '  -1.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.62'
$ws.Range('E30').Value = '  +13.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '37.63'
$ws.Range('E31').Value = '  +4.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '165.74'
$ws.Range('E32').Value = '  +1.41%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0887'
$ws.Range('E33').Value = '  +1.66%  '
$ws.Range('E34').Value = '  -3.88%  '
$ws.Range('E35').Value = '  -3.58%  '
$ws.Range('E36').Value = '  +0.98%  '
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('E38').Value = '  +1.77%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.75'
$ws.Range('E39').Value = '  +0.16%  '
$ws.Range('E40').Value = '  -1.14%  '
$ws.Range('E41').Value = '  +11.52%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '70.63'
$ws.Range('E42').Value = '  +2.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '96.32'
$ws.Range('E43').Value = '  -3.68%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.227'
$ws.Range('E44').Value = '  +0.62%  '
$ws.Range('E45').Value = '  +0.26%  '
$ws.Range('B46').Value = 'Celestia'
$ws.Range('C46').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.43'
$ws.Range('E46').Value = '  +3.35%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '117.55'
$ws.Range('E47').Value = '  +5.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '80.25'
$ws.Range('E48').Value = '  +4.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.654.32'
$ws.Range('E49').Value = '  +3.11%  '
$ws.Range('E50').Value = '  -0.22%  '
$ws.Range('E51').Value = '  +0.22%  '
